# Updates the crypto price/volume table on the active sheet to reflect the
# latest scrape (mirrors the GitHub Actions "Updated cryptos list" commit).
# Columns: D = Price (text), E = Volume(1h) (text, padded with spaces).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values like "1.00" / "67.658.56" as text. Force the
# whole column to Text format first so Excel doesn't reinterpret numeric-
# looking strings (e.g. "1.00" -> 1, "181.89" -> 181.89 as a number) when we
# assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.658.56"
$ws.Range("E2").Value = "  +4.55%  "

$ws.Range("D3").Value = "3.256.24"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "578.78"
$ws.Range("E5").Value = "  +2.55%  "

$ws.Range("D6").Value = "181.89"
$ws.Range("E6").Value = "  +6.57%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  -4.31%  "

$ws.Range("D9").Value = "3.257.24"
$ws.Range("E9").Value = "  +3.76%  "

$ws.Range("E10").Value = "  +4.86%  "

$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("E12").Value = "  +5.16%  "

$ws.Range("D13").Value = "3.828.58"
$ws.Range("E13").Value = "  +3.64%  "

$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("D15").Value = "28.14"
$ws.Range("E15").Value = "  +3.66%  "

$ws.Range("D16").Value = "67.633.34"
$ws.Range("E16").Value = "  +4.75%  "

$ws.Range("E17").Value = "  +3.14%  "

$ws.Range("D18").Value = "3.259.13"
$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("E19").Value = "  +1.97%  "

$ws.Range("D20").Value = "13.53"
$ws.Range("E20").Value = "  +5.06%  "

$ws.Range("D21").Value = "376.27"
$ws.Range("E21").Value = "  +5.92%  "

$ws.Range("E22").Value = "  +5.57%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").Value = "71.34"
$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("D25").Value = "0.512"
$ws.Range("E25").Value = "  +2.54%  "

$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("E28").Value = "  +3.46%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D30").Value = "5.74"
$ws.Range("E30").Value = "  +7.19%  "

$ws.Range("E31").Value = "  +4.25%  "

$ws.Range("E32").Value = "  +3.29%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("E34").Value = "  +6.16%  "

$ws.Range("D35").Value = "6.89"
$ws.Range("E35").Value = "  +3.75%  "

$ws.Range("D36").Value = "162.31"
$ws.Range("E36").Value = "  +5.74%  "

$ws.Range("E37").Value = "  +3.81%  "

$ws.Range("E38").Value = "  +3.29%  "

$ws.Range("D39").Value = "1.86"
$ws.Range("E39").Value = "  +7.67%  "

$ws.Range("D40").Value = "6.83"
$ws.Range("E40").Value = "  +13.28%  "

$ws.Range("D41").Value = "26.94"
$ws.Range("E41").Value = "  +3.48%  "

$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  +3.76%  "

$ws.Range("D43").Value = "363.87"
$ws.Range("E43").Value = "  +13.79%  "

$ws.Range("D44").Value = "4.46"
$ws.Range("E44").Value = "  +7.01%  "

$ws.Range("D45").Value = "2.749.43"
$ws.Range("E45").Value = "  +3.63%  "

$ws.Range("D46").Value = "25.65"
$ws.Range("E46").Value = "  +5.94%  "

$ws.Range("D47").Value = "40.65"
$ws.Range("E47").Value = "  +3.92%  "

$ws.Range("E48").Value = "  +3.30%  "

$ws.Range("E49").Value = "  +2.95%  "

$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +6.85%  "

$ws.Range("E51").Value = "  +0.22%  "
